$d = $word.ActiveDocument

$replacements = @(
    @("509÷7=72, 5", "810÷3=270, 0"),
    @("643÷3=214, 1", "657÷7=93, 6"),
    @("772÷3=257, 1", "151÷8=18, 7"),
    @("523÷6=87, 1", "524÷6=87, 2"),
    @("223÷8=27, 7", "724÷4=181, 0"),
    @("364÷6=60, 4", "921÷9=102, 3"),
    @("526÷7=75, 1", "324÷5=64, 4"),
    @("792÷7=113, 1", "674÷8=84, 2"),
    @("660÷9=73, 3", "223÷3=74, 1"),
    @("406÷6=67, 4", "604÷8=75, 4"),
    @("665÷3=221, 2", "624÷2=312, 0"),
    @("727÷4=181, 3", "108÷7=15, 3"),
    @("757÷4=189, 1", "305÷8=38, 1"),
    @("118÷6=19, 4", "610÷2=305, 0"),
    @("129÷4=32, 1", "635÷5=127, 0"),
    @("471÷8=58, 7", "891÷5=178, 1"),
    @("421÷9=46, 7", "974÷8=121, 6"),
    @("631÷7=90, 1", "205÷6=34, 1"),
    @("646÷4=161, 2", "359÷6=59, 5"),
    @("367÷4=91, 3", "826÷3=275, 1"),
    @("275÷9=30, 5", "362÷8=45, 2"),
    @("758÷3=252, 2", "396÷4=99, 0"),
    @("445÷2=222, 1", "369÷4=92, 1"),
    @("905÷8=113, 1", "204÷8=25, 4"),
    @("318÷5=63, 3", "599÷3=199, 2")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
